$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-26 21:48:42"
$ws.Range("O2").Value = "5.6 °C"
$ws.Range("E3").Value = "2026-02-26 21:48:45"
$ws.Range("L3").Value = "19.8 km/h - 114º 21:23 TU"
$ws.Range("E4").Value = "2026-02-26 21:48:48"
$ws.Range("H4").Value = "80%"
$ws.Range("O4").Value = "10.6 °C"
$ws.Range("E5").Value = "2026-02-26 21:48:50"
$ws.Range("E6").Value = "2026-02-26 21:48:53"
$ws.Range("E7").Value = "2026-02-26 21:48:56"
$ws.Range("E8").Value = "2026-02-26 21:48:58"
$ws.Range("J8").Value = "1027.0 hPa"
$ws.Range("E9").Value = "2026-02-26 21:49:01"
$ws.Range("O9").Value = "12.0 °C"
$ws.Range("E10").Value = "2026-02-26 21:49:04"
$ws.Range("H10").Value = "88%"
$ws.Range("K10").Value = "14.5 MJ/m2"
$ws.Range("O10").Value = "9.4 °C"
$ws.Range("E11").Value = "2026-02-26 21:49:06"
$ws.Range("O11").Value = "8.6 °C"
$ws.Range("E12").Value = "2026-02-26 21:49:09"
$ws.Range("E13").Value = "2026-02-26 21:49:11"
$ws.Range("O13").Value = "7.2 °C"
$ws.Range("E14").Value = "2026-02-26 21:49:14"
$ws.Range("H14").Value = "89%"
$ws.Range("O14").Value = "11.3 °C"
$ws.Range("E15").Value = "2026-02-26 21:49:16"
$ws.Range("O15").Value = "11.5 °C"
$ws.Range("E16").Value = "2026-02-26 21:49:19"
$ws.Range("E17").Value = "2026-02-26 21:49:22"
$ws.Range("H17").Value = "37%"
$ws.Range("O17").Value = "8.4 °C"
$ws.Range("E18").Value = "2026-02-26 21:49:24"
$ws.Range("E19").Value = "2026-02-26 21:49:27"
$ws.Range("H19").Value = "48%"
$ws.Range("E20").Value = "2026-02-26 21:49:30"
$ws.Range("E21").Value = "2026-02-26 21:49:32"
$ws.Range("E22").Value = "2026-02-26 21:49:35"
$ws.Range("E23").Value = "2026-02-26 21:49:38"
$ws.Range("H23").Value = "39%"
$ws.Range("E24").Value = "2026-02-26 21:49:41"
$ws.Range("E25").Value = "2026-02-26 21:49:43"
$ws.Range("E26").Value = "2026-02-26 21:49:46"
$ws.Range("H26").Value = "42%"
$ws.Range("K26").Value = "16.3 MJ/m2"
$ws.Range("E27").Value = "2026-02-26 21:49:48"
$ws.Range("L27").Value = "20.9 km/h - 242º 21:17 TU"
$ws.Range("E28").Value = "2026-02-26 21:49:51"
$ws.Range("H28").Value = "79%"
$ws.Range("N28").Value = "5.7 °C 21:14 TU"
$ws.Range("O28").Value = "10.7 °C"
$ws.Range("E29").Value = "2026-02-26 21:49:54"
$ws.Range("N29").Value = "7.1 °C 21:03 TU"
$ws.Range("O29").Value = "11.5 °C"
$ws.Range("E30").Value = "2026-02-26 21:49:57"
$ws.Range("H30").Value = "87%"
$ws.Range("O30").Value = "12.0 °C"
$ws.Range("E31").Value = "2026-02-26 21:49:59"
$ws.Range("E32").Value = "2026-02-26 21:50:02"
$ws.Range("H32").Value = "67%"
$ws.Range("O32").Value = "7.9 °C"
$ws.Range("E33").Value = "2026-02-26 21:50:05"
$ws.Range("J33").Value = "1026.8 hPa"
$ws.Range("O33").Value = "8.6 °C"
$ws.Range("E34").Value = "2026-02-26 21:50:07"
$ws.Range("O34").Value = "4.7 °C"
$ws.Range("E35").Value = "2026-02-26 21:50:10"
$ws.Range("E36").Value = "2026-02-26 21:50:13"
$ws.Range("E37").Value = "2026-02-26 21:50:15"
$ws.Range("J37").Value = "1028.3 hPa"
$ws.Range("O37").Value = "7.7 °C"
$ws.Range("E38").Value = "2026-02-26 21:50:18"
$ws.Range("H38").Value = "81%"
$ws.Range("O38").Value = "11.1 °C"
$ws.Range("E39").Value = "2026-02-26 21:50:21"
$ws.Range("E40").Value = "2026-02-26 21:50:23"
$ws.Range("O40").Value = "9.5 °C"
$ws.Range("E41").Value = "2026-02-26 21:50:26"
$ws.Range("K41").Value = "15.9 MJ/m2"
$ws.Range("E42").Value = "2026-02-26 21:50:29"
$ws.Range("N42").Value = "7.3 °C 21:20 TU"
$ws.Range("O42").Value = "11.2 °C"
$ws.Range("E43").Value = "2026-02-26 21:50:31"
$ws.Range("E44").Value = "2026-02-26 21:50:34"
$ws.Range("O44").Value = "2.2 °C"
$ws.Range("E45").Value = "2026-02-26 21:50:37"
$ws.Range("O45").Value = "10.5 °C"
$ws.Range("E46").Value = "2026-02-26 21:50:40"
